# Insert two new price rows at the top of the Femacal de La Calera - Plátano
# data block (rows 984-985), shifting all subsequent rows down by two.
# The rest of the block (previously rows 984-1050) slides down to 986-1052
# unchanged; two brand-new rows carry the latest (2022-09-xx) price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("984:985").Insert()

# New row 984: Pintón
$ws.Range("A984").Value = 3
$ws.Range("B984").Value = "Femacal de La Calera"
$ws.Range("C984").Value = "Coquimbo"
$ws.Range("D984").Value = 44826
$ws.Range("E984").Value = 5
$ws.Range("F984").Value = "Fruta"
$ws.Range("G984").Value = 100108
$ws.Range("H984").Value = "Tropicales y subtropicales"
$ws.Range("I984").Value = 100108006
$ws.Range("J984").Value = "Plátano"
$ws.Range("K984").Value = "Sin especificar"
$ws.Range("L984").Value = "Pintón"
$ws.Range("M984").Value = 400
$ws.Range("N984").Value = 22000
$ws.Range("O984").Value = 22500
$ws.Range("P984").Value = 22200
$ws.Range("Q984").Value = "`$/caja 20 kilos"
$ws.Range("R984").Value = "Ecuador"
$ws.Range("S984").Value = 1110
$ws.Range("T984").Value = 20

# New row 985: Primera Pintón
$ws.Range("A985").Value = 3
$ws.Range("B985").Value = "Femacal de La Calera"
$ws.Range("C985").Value = "Coquimbo"
$ws.Range("D985").Value = 44826
$ws.Range("E985").Value = 5
$ws.Range("F985").Value = "Fruta"
$ws.Range("G985").Value = 100108
$ws.Range("H985").Value = "Tropicales y subtropicales"
$ws.Range("I985").Value = 100108006
$ws.Range("J985").Value = "Plátano"
$ws.Range("K985").Value = "Sin especificar"
$ws.Range("L985").Value = "Primera Pintón"
$ws.Range("M985").Value = 480
$ws.Range("N985").Value = 23000
$ws.Range("O985").Value = 24000
$ws.Range("P985").Value = 23583
$ws.Range("Q985").Value = "`$/caja 20 kilos"
$ws.Range("R985").Value = "Ecuador"
$ws.Range("S985").Value = 1179
$ws.Range("T985").Value = 20
